# Update the cryptos list with latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($cellRef, $value) {
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $value
}

Set-Text "D2" "41.103.90"
Set-Text "E2" "  -2.35%  "

Set-Text "D3" "2.179.53"
Set-Text "E3" "  -2.17%  "

Set-Text "E4" "  -0.04%  "

Set-Text "D5" "236.98"
Set-Text "E5" "  -2.94%  "

Set-Text "D6" "0.613"
Set-Text "E6" "  -1.45%  "

Set-Text "D7" "70.37"
Set-Text "E7" "  -4.67%  "

Set-Text "E8" "  +0.03%  "

Set-Text "D9" "0.578"
Set-Text "E9" "  -6.13%  "

Set-Text "D10" "40.19"
Set-Text "E10" "  -7.02%  "

Set-Text "D11" "0.0929"
Set-Text "E11" "  -3.90%  "

Set-Text "E12" "  -1.46%  "

Set-Text "D13" "6.78"
Set-Text "E13" "  -5.26%  "

Set-Text "D14" "2.504.79"
Set-Text "E14" "  -2.26%  "

Set-Text "D15" "13.96"
Set-Text "E15" "  -2.50%  "

Set-Text "D16" "0.812"
Set-Text "E16" "  -4.22%  "

Set-Text "D17" "2.191.99"
Set-Text "E17" "  -1.55%  "

Set-Text "D18" "41.029.77"
Set-Text "E18" "  -2.46%  "

Set-Text "D19" "0.0000102"
Set-Text "E19" "  -8.82%  "

Set-Text "D20" "70.62"
Set-Text "E20" "  -2.55%  "

Set-Text "D21" "5.97"
Set-Text "E21" "  -4.22%  "

Set-Text "D22" "9.86"
Set-Text "E22" "  -4.69%  "

Set-Text "D23" "226.33"
Set-Text "E23" "  -1.99%  "

Set-Text "D24" "1.97"
Set-Text "E24" "  -7.82%  "

Set-Text "E25" "  +0.04%  "

Set-Text "D26" "10.97"
Set-Text "E26" "  -6.66%  "

Set-Text "D27" "3.56"
Set-Text "E27" "  -0.84%  "

Set-Text "D28" "2.21"
Set-Text "E28" "  -3.51%  "

Set-Text "E29" "  -1.64%  "

Set-Text "D30" "166.38"
Set-Text "E30" "  -0.42%  "

Set-Text "D31" "19.95"
Set-Text "E31" "  -3.68%  "

# Rows 32 and 33: contents swapped (Hedera/InjectiveProtocol) with new values.
Set-Text "B32" "Hedera"
Set-Text "C32" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-Text "D32" "0.0780"
Set-Text "E32" "  -2.90%  "

Set-Text "B33" "InjectiveProtocol"
Set-Text "C33" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-Text "D33" "30.23"
Set-Text "E33" "  +2.12%  "

Set-Text "D34" "5.14"
Set-Text "E34" "  -8.60%  "

Set-Text "E35" "  -3.26%  "

Set-Text "D36" "0.105"
Set-Text "E36" "  -9.00%  "

Set-Text "D37" "4.15"
Set-Text "E37" "  -5.23%  "

Set-Text "E38" "  -4.58%  "

Set-Text "D39" "12.45"
Set-Text "E39" "  -5.46%  "

Set-Text "D40" "2.06"
Set-Text "E40" "  -4.45%  "

Set-Text "D41" "5.45"
Set-Text "E41" "  -3.22%  "

Set-Text "D42" "60.35"
Set-Text "E42" "  -6.95%  "

Set-Text "D43" "0.191"
Set-Text "E43" "  -5.43%  "

Set-Text "D44" "8.35"
Set-Text "E44" "  -4.65%  "

Set-Text "D45" "0.0978"
Set-Text "E45" "  -3.68%  "

Set-Text "D46" "99.29"
Set-Text "E46" "  -5.90%  "

Set-Text "D47" "1.09"
Set-Text "E47" "  -2.69%  "

Set-Text "D48" "1.13"
Set-Text "E48" "  -3.44%  "

Set-Text "D49" "2.23"
Set-Text "E49" "  -8.20%  "

Set-Text "E50" "  -1.98%  "

Set-Text "D51" "2.381.37"
Set-Text "E51" "  -2.28%  "
